# Update WiFly Message IDs sheet with the new message ID table.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New data table (rows 2-13), replacing the previous IRxy / ACK rows.
$data = @(
    @("IR0",   "0x10"),
    @("IR1",   "0x11"),
    @("IR2",   "0x12"),
    @("IR3",   "0x13"),
    @("IR4",   "0x14"),
    @("IR5",   "0x15"),
    @("IR6",   "0x16"),
    @("IR7",   "0x17"),
    @("IR8",   "0x18"),
    @("IR9",   "0x19"),
    @("MOTOR", "0x50"),
    @("ACK",   "0x99")
)

for ($i = 0; $i -lt $data.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 1).Value = $data[$i][0]
    $ws.Cells.Item($row, 2).Value = $data[$i][1]
}

# Column A needs to widen to fit the new longer "Message Type" / values.
$ws.Columns.Item(1).ColumnWidth = 20.33

# Update the active selection to match the post-edit state.
$ws.Range("C12").Select()
